$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 49.96455633333333
$ws.Range("H2").Value = 149.893669
$ws.Range("I2").Value = 0.1551859508057627
$ws.Range("J2").Value = 0.1551859508057627
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 189.6251011346681
$ws.Range("R2").Value = 1706.625910212013
$ws.Range("S2").Value = 0.001620473045247292
$ws.Range("T2").Value = 0.001620473045247292
$ws.Range("G3").Value = 49.96455633333333
$ws.Range("H3").Value = 149.893669
$ws.Range("I3").Value = 0.1551859508057627
$ws.Range("J3").Value = 0.1551859508057627
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 12160.19378138446
$ws.Range("R3").Value = 109441.7440324601
$ws.Range("S3").Value = 0.1039169715918725
$ws.Range("T3").Value = 0.1039169715918725
$ws.Range("G4").Value = 49.96455633333333
$ws.Range("H4").Value = 149.893669
$ws.Range("I4").Value = 0.1551859508057627
$ws.Range("J4").Value = 0.1551859508057627
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 1489.012280140066
$ws.Range("R4").Value = 13401.11052126059
$ws.Range("S4").Value = 0.01272460370262684
$ws.Range("T4").Value = 0.01272460370262684
$ws.Range("G5").Value = 49.96455633333333
$ws.Range("H5").Value = 149.893669
$ws.Range("I5").Value = 0.1551859508057627
$ws.Range("J5").Value = 0.1551859508057627
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 4320.774578719653
$ws.Range("R5").Value = 38886.97120847688
$ws.Range("S5").Value = 0.03692390246601609
$ws.Range("T5").Value = 0.03692390246601609
$ws.Range("I6").Value = 0.5303393919600503
$ws.Range("J6").Value = 0.5303393919600503
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 648.0332808090029
$ws.Range("R6").Value = 5832.299527281025
$ws.Range("S6").Value = 0.00553787688280985
$ws.Range("T6").Value = 0.00553787688280985
$ws.Range("I7").Value = 0.5303393919600503
$ws.Range("J7").Value = 0.5303393919600503
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("S7").Value = 0.355130495010744
$ws.Range("T7").Value = 0.355130495010744
$ws.Range("I8").Value = 0.5303393919600503
$ws.Range("J8").Value = 0.5303393919600503
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 5088.616998963585
$ws.Range("R8").Value = 45797.55299067226
$ws.Range("S8").Value = 0.04348562840608074
$ws.Range("T8").Value = 0.04348562840608073
$ws.Range("I9").Value = 0.5303393919600503
$ws.Range("J9").Value = 0.5303393919600503
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 14766.00781821244
$ws.Range("R9").Value = 132894.0703639119
$ws.Range("S9").Value = 0.1261853916604157
$ws.Range("T9").Value = 0.1261853916604157
$ws.Range("G10").Value = 13.36927633333333
$ws.Range("H10").Value = 40.107829
$ws.Range("I10").Value = 0.04152391238164931
$ws.Range("J10").Value = 0.04152391238164931
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 50.73897504248145
$ws.Range("R10").Value = 456.650775382333
$ws.Range("S10").Value = 0.0004335984049992642
$ws.Range("T10").Value = 0.0004335984049992642
$ws.Range("G11").Value = 13.36927633333333
$ws.Range("H11").Value = 40.107829
$ws.Range("I11").Value = 0.04152391238164931
$ws.Range("J11").Value = 0.04152391238164931
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 3253.766326786165
$ws.Range("R11").Value = 29283.89694107548
$ws.Range("S11").Value = 0.02780560483047939
$ws.Range("T11").Value = 0.02780560483047939
$ws.Range("G12").Value = 13.36927633333333
$ws.Range("H12").Value = 40.107829
$ws.Range("I12").Value = 0.04152391238164931
$ws.Range("J12").Value = 0.04152391238164931
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 398.4227640111863
$ws.Range("R12").Value = 3585.804876100676
$ws.Range("S12").Value = 0.003404788426372591
$ws.Range("T12").Value = 0.003404788426372591
$ws.Range("G13").Value = 13.36927633333333
$ws.Range("H13").Value = 40.107829
$ws.Range("I13").Value = 0.04152391238164931
$ws.Range("J13").Value = 0.04152391238164931
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 1156.13213758104
$ws.Range("R13").Value = 10405.18923822936
$ws.Range("S13").Value = 0.009879920719798055
$ws.Range("T13").Value = 0.009879920719798055
$ws.Range("G14").Value = 87.88078300000001
$ws.Range("H14").Value = 263.642349
$ws.Range("I14").Value = 0.2729507448525377
$ws.Range("J14").Value = 0.2729507448525377
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 333.5244738889304
$ws.Range("R14").Value = 3001.720265000374
$ws.Range("S14").Value = 0.002850189224070427
$ws.Range("T14").Value = 0.002850189224070427
$ws.Range("G15").Value = 87.88078300000001
$ws.Range("H15").Value = 263.642349
$ws.Range("I15").Value = 0.2729507448525377
$ws.Range("J15").Value = 0.2729507448525377
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 21388.10847854683
$ws.Range("R15").Value = 192492.9763069215
$ws.Range("S15").Value = 0.1827756614020005
$ws.Range("T15").Value = 0.1827756614020005
$ws.Range("G16").Value = 87.88078300000001
$ws.Range("H16").Value = 263.642349
$ws.Range("I16").Value = 0.2729507448525377
$ws.Range("J16").Value = 0.2729507448525377
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 2618.967817953493
$ws.Range("R16").Value = 23570.71036158144
$ws.Range("S16").Value = 0.0223808279071122
$ws.Range("T16").Value = 0.02238082790711219
$ws.Range("G17").Value = 87.88078300000001
$ws.Range("H17").Value = 263.642349
$ws.Range("I17").Value = 0.2729507448525377
$ws.Range("J17").Value = 0.2729507448525377
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 7599.648250875325
$ws.Range("R17").Value = 68396.83425787793
$ws.Range("S17").Value = 0.06494406631935452
$ws.Range("T17").Value = 0.06494406631935452
